# Commit: made type split a little more complex
# Column B ("pre-1955" vintage stock) on Sheet1 is re-baselined: rows for
# years 1600-1900 (B2:B302) drop to 0, and the remaining years (1901-2050,
# B303:B452) are replaced with a new, smaller stock-split curve. All other
# columns (A, C:H) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newB = @(
    0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 
    0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 
    0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 
    0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 
    0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 
    0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 
    0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 
    0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 
    0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 
    0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 148821.280951146, 299955.2009834384, 453440.9515181568, 
    609318.2193980068, 767627.1745164256, 928408.4558839672, 1091703.156070428, 1257552.803961648, 
    1425999.34577052, 1597085.124242926, 1770852.856000254, 1947345.60696158, 2126606.765790498, 
    2308680.015313555, 2493609.301859512, 2681438.802471632, 2872212.889947923, 3065976.09566791, 
    3262773.070168414, 3462648.541434649, 3665647.270877997, 3871814.00697635, 4081193.43655888, 
    4293830.133722389, 4509768.506372998, 4729052.740393559, 4951726.741444175, 5177834.074410799, 
    5407417.900524851, 5640520.912185107, 5877185.265521674, 6117452.510751228, 6361363.520381849, 
    6608958.41533574, 6860276.489068064, 7115356.129770517, 7374234.740758609, 7636948.65915283, 
    7903533.072974239, 8174021.936786264, 8448447.886025677, 8726842.150176479, 9009234.464951754, 
    9295652.983658932, 9586124.18793539, 9880672.798050757, 10179321.68298355, 10482091.77048836, 
    10789001.95738027, 11100069.02027082, 11415307.52699907, 11734729.74900785, 12058345.57492265, 
    12386162.42559634, 12718185.17088768, 12715661.31998215, 12712960.68505145, 12710072.32312029, 
    12706984.72098878, 12703685.773024, 12700162.75860848, 12696402.31928344, 12692390.43562913, 
    12688112.40392988, 12683552.81267641, 12678695.51896385, 12673523.62484977, 12668019.45374267, 
    12662164.52689846, 12655939.54010892, 12649324.34067396, 12642297.90475696, 12634838.31523054, 
    12626922.74012859, 12618527.41182882, 12609627.60709943, 12600197.62815229, 12590210.78485458, 
    12579639.37826073, 12568454.6856357, 12556626.94715073, 12544125.35444245, 12530918.04123586, 
    12516972.07624109, 12502253.45854344, 12486727.11571476, 12470356.90488277, 12453105.61700279, 
    12434934.98458341, 12415805.69312375, 12395677.39652522, 12374508.73674475, 12352257.367959, 
    12328879.98551037, 12304332.35990492, 12278569.37613016, 12251545.0785561, 12223212.72167637, 
    12193524.82693726, 12162433.24589104, 12129889.22989558, 12095843.50656524, 12060246.36315813, 
    12023047.73706135, 11984197.31350982, 11943644.63064441, 11901339.19198233, 11857230.58633651, 
    11811268.61518084, 11763403.42741553, 11713585.66144077, 11661766.59439742, 11607898.29838188, 
    11551933.80338699, 11493827.26666412, 11433534.14814179, 11371011.39147518, 11306217.61023838, 
    11239113.27870739, 11169660.92661786, 11097825.33721763, 11023573.74787029, 10946876.05240372, 
    10867705.00433721, 10786036.42006192, 10701849.38099528, 10615126.43367766, 10525853.78673356, 
    10434021.50357744, 10339623.6897086, 10242658.67341016, 10143129.17864506, 10041042.48892728, 
    9936410.600940641, 9829250.366679834, 9719583.622900531, 9607437.306686671, 9492843.555974679, 
    9375839.793915791, 9256468.796009673, 9134778.739004539, 9010823.230631426, 8884661.31932251, 
    8756357.483155785, 8625981.597369405, 8493608.879899669, 8359319.814514743, 8223200.051242239, 
    8085340.28392117, 7945836.104846989 
)

$arr = New-Object 'object[,]' $newB.Length,1
for ($i = 0; $i -lt $newB.Length; $i++) {
    $arr[$i,0] = $newB[$i]
}

$ws.Range("B2:B452").Value = $arr
